$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: the "Serienummer" value moves from "1908220400048819" to "1908222340013186"
# (the value that used to live in row 4), and the stray styled-but-empty W3 cell is cleared.
$ws.Range("V3").Value = "1908222340013186"
$ws.Range("W3").Clear()

# Rows 4 and 5 are removed entirely (full row delete, shifting rows up).
$ws.Rows("4:5").Delete()

# Update the active selection to match the new state.
$ws.Range("W4").Select()

# Restore the window position recorded in the workbook view.
$excel.ActiveWindow.Left = 5025
$excel.ActiveWindow.Top = 1455
